# Add two new weekly survey columns (BJ/BK on "data", BI/BJ on "pocetR")
# covering the weeks "9.-15. 8. 2021" and "16.-22. 8. 2021", and refresh the
# "aktualizace" (updated-on) date stamps from 27. 7. 2021 to 1. 9. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data" (first sheet): new columns BJ (col 62) and BK (col 63)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)

# Header row: copy formatting from the previous header cell (BI1) so the new
# header cells pick up the same style (bold, centered, bordered), then set
# the correct text.
$wsData.Range("BI1").Copy($wsData.Range("BJ1"))
$wsData.Range("BI1").Copy($wsData.Range("BK1"))
$wsData.Range("BJ1").Value = "9.–15. 8. 2021"
$wsData.Range("BK1").Value = "16.–22. 8. 2021"

# Data rows 2-77: two new numeric values per row (BJ, BK).
$dataRows = @(
    @(2, 0.64, 0.63),
    @(3, 0.07000000000000001, 0.07000000000000001),
    @(4, 0.07000000000000001, 0.06),
    @(5, 0.22, 0.24),
    @(6, 0.67, 0.6),
    @(7, 0.04, 0.05),
    @(8, 0.08, 0.07000000000000001),
    @(9, 0.21, 0.28),
    @(10, 0.42, 0.4),
    @(11, 0.19, 0.15),
    @(12, 0.25, 0.22),
    @(13, 0.14, 0.23),
    @(14, 0.74, 0.75),
    @(15, 0.05, 0.08),
    @(16, 0.02, 0.015),
    @(17, 0.19, 0.155),
    @(18, 0.71, 0.71),
    @(19, 0.05, 0.03),
    @(20, 0.02, 0.03),
    @(21, 0.22, 0.23),
    @(22, 0.31, 0.38),
    @(23, 0.11, 0.12),
    @(24, 0.04, 0.02),
    @(25, 0.54, 0.48),
    @(26, 0.66, 0.65),
    @(27, 0.06, 0.06),
    @(28, 0.07000000000000001, 0.05),
    @(29, 0.21, 0.24),
    @(30, 0.71, 0.67),
    @(31, 0.04, 0.05),
    @(32, 0.04, 0.04),
    @(33, 0.21, 0.24),
    @(34, 0.5600000000000001, 0.52),
    @(35, 0.08, 0.12),
    @(36, 0.14, 0.13),
    @(37, 0.22, 0.23),
    @(38, 0.53, 0.5600000000000001),
    @(39, 0.12, 0.11),
    @(40, 0.11, 0.09),
    @(41, 0.24, 0.24),
    @(42, 0.64, 0.62),
    @(43, 0.08, 0.08),
    @(44, 0.11, 0.11),
    @(45, 0.17, 0.19),
    @(46, 0.67, 0.65),
    @(47, 0.06, 0.07000000000000001),
    @(48, 0.05, 0.04),
    @(49, 0.22, 0.24),
    @(50, 0.55, 0.54),
    @(51, 0.08, 0.08),
    @(52, 0.08, 0.06),
    @(53, 0.29, 0.32),
    @(54, 0.62, 0.64),
    @(55, 0.08, 0.07000000000000001),
    @(56, 0.07000000000000001, 0.05),
    @(57, 0.23, 0.24),
    @(58, 0.58, 0.5600000000000001),
    @(59, 0.09, 0.11),
    @(60, 0.05, 0.03),
    @(61, 0.28, 0.3),
    @(62, 0.66, 0.64),
    @(63, 0.06, 0.06),
    @(64, 0.08, 0.07000000000000001),
    @(65, 0.2, 0.23),
    @(66, 0.75, 0.72),
    @(67, 0.02, 0.03),
    @(68, 0.04, 0.02),
    @(69, 0.19, 0.23),
    @(70, 0.62, 0.61),
    @(71, 0.08, 0.09),
    @(72, 0.09, 0.07000000000000001),
    @(73, 0.21, 0.23),
    @(74, 0.51, 0.51),
    @(75, 0.11, 0.1),
    @(76, 0.1, 0.11),
    @(77, 0.28, 0.28)
)
foreach ($r in $dataRows) {
    $wsData.Cells.Item($r[0], 62).Value = $r[1]
    $wsData.Cells.Item($r[0], 63).Value = $r[2]
}

# Row 78 footer label: bump the "aktualizace" date.
$wsData.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" (second sheet): new columns BI (col 61) and BJ (col 62)
# ---------------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item(2)

$wsPocet.Range("BH1").Copy($wsPocet.Range("BI1"))
$wsPocet.Range("BH1").Copy($wsPocet.Range("BJ1"))
$wsPocet.Range("BI1").Value = "9.–15. 8. 2021"
$wsPocet.Range("BJ1").Value = "16.–22. 8. 2021"

$pocetRows = @(
    @(2, 1073, 1073),
    @(3, 293, 293),
    @(4, 91, 91),
    @(5, 292, 292),
    @(6, 158, 158),
    @(7, 85, 85),
    @(8, 523, 523),
    @(9, 262, 262),
    @(10, 136, 136),
    @(11, 152, 152),
    @(12, 267, 267),
    @(13, 620, 620),
    @(14, 186, 186),
    @(15, 211, 211),
    @(16, 192, 192),
    @(17, 670, 670),
    @(18, 387, 387),
    @(19, 436, 436),
    @(20, 250, 250)
)
foreach ($r in $pocetRows) {
    $wsPocet.Cells.Item($r[0], 61).Value = $r[1]
    $wsPocet.Cells.Item($r[0], 62).Value = $r[2]
}

# Row 21 footer label + two trailing blank placeholder cells (BI21, BJ21)
# matching the other empty cells in that row.
$wsPocet.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"
$wsPocet.Range("B21").Copy($wsPocet.Range("BI21"))
$wsPocet.Range("B21").Copy($wsPocet.Range("BJ21"))
